# SAV-700: Update charts importer test fixture
# - Ensure ID of all complex chart core questions is constant
#   (Core sheet: code column (A) now mirrors the type column (B) instead of
#   using ad-hoc testchartcorecode* identifiers)
# - Ensure ID of charting date recorded program data element is constant
#   (Test Chart sheet: the charting date question code becomes the constant
#   "PatientChartingDate" instead of "testchartcode0")

$wb = $excel.ActiveWorkbook

$core = $wb.Worksheets.Item("Core")
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartSubtype"
$core.Range("A5").Value = "ComplexChartType"

$testChart = $wb.Worksheets.Item("Test Chart")
$testChart.Range("A2").Value = "PatientChartingDate"
$testChart.Range("A2").HorizontalAlignment = -4131
